# Auto-generated edit script applying numeric updates described in the diff
# for Sheets/Atomos_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3412
$ws.Range("J64").Value = 3334.6667
$ws.Range("L64").Value = 3334.6667
$ws.Range("N64").Value = -3830.6667

$ws.Range("H67").Value = 3412
$ws.Range("J67").Value = 3334.6667
$ws.Range("L67").Value = 3334.6667
$ws.Range("N67").Value = -5050.6667

$ws.Range("H76").Value = 3584.85
$ws.Range("I76").Value = 3299.4285
$ws.Range("J76").Value = 4250.8335
$ws.Range("K76").Value = 3299.4285
$ws.Range("L76").Value = 4250.8335
$ws.Range("M76").Value = -2984.4285
$ws.Range("N76").Value = -4880.8335

$ws.Range("H79").Value = 3584.85
$ws.Range("I79").Value = 3299.4285
$ws.Range("J79").Value = 4250.8335
$ws.Range("K79").Value = 3299.4285
$ws.Range("L79").Value = 4250.8335
$ws.Range("M79").Value = -2207.4285
$ws.Range("N79").Value = -6434.8335

$ws.Range("H116").Value = 3067.9363
$ws.Range("I116").Value = 2247.1614
$ws.Range("J116").Value = 4658.1875
$ws.Range("K116").Value = 2247.1614
$ws.Range("L116").Value = 4658.1875
$ws.Range("M116").Value = 1194.8386
$ws.Range("N116").Value = -11542.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35716572
$ws.Range("I2").Value = 35716572
$ws.Range("K2").Value = 35716572
$ws.Range("M2").Value = -35716459

$ws.Range("H37").Value = 23394
$ws.Range("J37").Value = 45038
$ws.Range("L37").Value = 45038
$ws.Range("N37").Value = -45584

$ws.Range("H74").Value = 1479
$ws.Range("I74").Value = 1432.3334
$ws.Range("K74").Value = 1432.3334
$ws.Range("M74").Value = -558.3334

$ws.Range("H77").Value = 1479
$ws.Range("I77").Value = 1432.3334
$ws.Range("K77").Value = 7161.666999999999
$ws.Range("M77").Value = -2793.666999999999

$ws.Range("H88").Value = 2331.2
$ws.Range("I88").Value = 1968.6666
$ws.Range("K88").Value = 1968.6666
$ws.Range("M88").Value = -1562.6666

$ws.Range("H91").Value = 2331.2
$ws.Range("I91").Value = 1968.6666
$ws.Range("K91").Value = 1968.6666
$ws.Range("M91").Value = -564.6666

$ws.Range("H116").Value = 35716572
$ws.Range("I116").Value = 35716572
$ws.Range("K116").Value = 35716572
$ws.Range("M116").Value = -35714278

$ws.Range("H132").Value = 28575420
$ws.Range("I132").Value = 35717916
$ws.Range("J132").Value = 5429.857
$ws.Range("K132").Value = 107153748
$ws.Range("L132").Value = 16289.571
$ws.Range("M132").Value = -107151218
$ws.Range("N132").Value = -21349.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35716572
$ws.Range("I3").Value = 35716572
$ws.Range("K3").Value = 35716572
$ws.Range("M3").Value = -35716458

$ws.Range("H86").Value = 1760.1305
$ws.Range("I86").Value = 1431
$ws.Range("J86").Value = 2512.4285
$ws.Range("K86").Value = 1431
$ws.Range("L86").Value = 2512.4285
$ws.Range("M86").Value = -308
$ws.Range("N86").Value = -4758.4285

$ws.Range("H89").Value = 1760.1305
$ws.Range("I89").Value = 1431
$ws.Range("J89").Value = 2512.4285
$ws.Range("K89").Value = 7155
$ws.Range("L89").Value = 12562.1425
$ws.Range("M89").Value = -1539
$ws.Range("N89").Value = -23794.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1589107.1
$ws.Range("I31").Value = 1696165.2
$ws.Range("K31").Value = 1696165.2
$ws.Range("M31").Value = -1695870.2

$ws.Range("H34").Value = 1589107.1
$ws.Range("I34").Value = 1696165.2
$ws.Range("K34").Value = 1696165.2
$ws.Range("M34").Value = -1695963.2

$ws.Range("H105").Value = 2033.56
$ws.Range("I105").Value = 1697.2273
$ws.Range("K105").Value = 1697.2273
$ws.Range("M105").Value = 49.77269999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1946.8572
$ws.Range("I5").Value = 744.6
$ws.Range("J5").Value = 4952.5
$ws.Range("K5").Value = 2233.8
$ws.Range("L5").Value = 14857.5
$ws.Range("M5").Value = -2121.8
$ws.Range("N5").Value = -15081.5

$ws.Range("H122").Value = 1889
$ws.Range("J122").Value = 2043.3334
$ws.Range("L122").Value = 18390.0006
$ws.Range("N122").Value = -23290.0006

$ws.Range("H131").Value = 1137.5264
$ws.Range("J131").Value = 1083.7858
$ws.Range("L131").Value = 3251.3574
$ws.Range("N131").Value = -13331.3574

$ws.Range("H135").Value = 1946.8572
$ws.Range("I135").Value = 744.6
$ws.Range("J135").Value = 4952.5
$ws.Range("K135").Value = 6701.400000000001
$ws.Range("L135").Value = 44572.5
$ws.Range("M135").Value = -4166.400000000001
$ws.Range("N135").Value = -49642.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5008
$ws.Range("I70").Value = 5008
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5008
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4738
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 5008
$ws.Range("I73").Value = 5008
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5008
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4072
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 4875
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 5400
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 5400
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -7396

$ws.Range("H83").Value = 4875
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 5400
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 27000
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -36984

$ws.Range("H102").Value = 788069.0600000001
$ws.Range("I102").Value = 1996.1875
$ws.Range("J102").Value = 4980457.5
$ws.Range("K102").Value = 1996.1875
$ws.Range("L102").Value = 4980457.5
$ws.Range("M102").Value = -374.1875
$ws.Range("N102").Value = -4983701.5

$ws.Range("H132").Value = 2625.3438
$ws.Range("I132").Value = 1765.7826
$ws.Range("K132").Value = 5297.3478
$ws.Range("M132").Value = -2767.3478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1365.2
$ws.Range("I55").Value = 216.66667
$ws.Range("J55").Value = 1857.4286
$ws.Range("K55").Value = 216.66667
$ws.Range("L55").Value = 1857.4286
$ws.Range("M55").Value = -43.66667000000001
$ws.Range("N55").Value = -2203.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4427487.5
$ws.Range("I126").Value = 1543.85
$ws.Range("K126").Value = 4631.549999999999
$ws.Range("M126").Value = -2161.549999999999
